$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the sample data rows (2-7), leaving only the header row.
# Row 4's I4 cell keeps its existing style (bold "drop-flag" highlight)
# even though its value is cleared, so that row survives as a blank
# template row; rows 2, 3, 5, 6, 7 become fully empty and are dropped.
$ws.Range("A2:J7").ClearContents()

# Match the author's new selection/active cell left behind in the file.
$ws.Range("A2:I7").Select()
